# Applies the cryptos-list data refresh described in the commit:
#   "Updated cryptos list on Fri Jun 16 11:08:49 UTC 2023 with GitHub Actions"
#
# The source sheet stores every data cell as plain text (inline strings),
# including numeric-looking prices such as "0.9989" or "1.0000". If we just
# assign .Value with such a string, Excel COM helpfully (and wrongly, for
# this sheet) reinterprets it as a number -- losing trailing zeros / exact
# formatting (e.g. "1.0000" -> 1, "0.06160" -> 0.0616). To avoid that we
# force each touched cell's number format to Text ("@") before writing it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cell address -> new value, in the order the diff applies them
$updates = [ordered]@{
    D2 = '25.542.05'
    E2 = '  +2.71%  '
    D3 = '1.670.20'
    E3 = '  +2.29%  '
    D4 = '0.9989'
    E4 = '  +0.06%  '
    D5 = '235.60'
    D7 = '0.4738'
    E7 = '  +0.53%  '
    D8 = '0.2598'
    E8 = '  +1.84%  '
    D9 = '0.06160'
    E9 = '  +1.55%  '
    D10 = '1.667.07'
    E10 = '  +2.09%  '
    D11 = '0.06995'
    E11 = '  +0.50%  '
    D12 = '14.77'
    E12 = '  +1.74%  '
    D13 = '0.5846'
    E13 = '  -1.85%  '
    D14 = '4.367'
    E14 = '  +1.47%  '
    D15 = '75.45'
    E15 = '  +3.17%  '
    D17 = '0.9994'
    E17 = '  +0.04%  '
    D18 = '25.544.09'
    E18 = '  +2.70%  '
    D19 = '0.000006717'
    E19 = '  +2.71%  '
    D20 = '11.40'
    E20 = '  +2.73%  '
    D21 = '1.881.68'
    E21 = '  +1.68%  '
    D22 = '4.439'
    E22 = '  +2.80%  '
    D23 = '8.757'
    E23 = '  +2.70%  '
    D24 = '5.218'
    E24 = '  +0.23%  '
    D25 = '136.71'
    E25 = '  +2.82%  '
    D26 = '14.96'
    E26 = '  +1.55%  '
    D27 = '1.386'
    E27 = '  +0.82%  '
    D28 = '1.714'
    E28 = '  +5.86%  '
    D29 = '104.36'
    E29 = '  +1.03%  '
    E30 = '  +5.10%  '
    D31 = '0.07831'
    E31 = '  +1.86%  '
    E32 = '  +3.23%  '
    B33 = 'Hedera'
    C33 = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
    D33 = '0.04294'
    E33 = '  +0.46%  '
    B34 = 'HuobiToken'
    C34 = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
    D34 = '2.624'
    E34 = '  +1.70%  '
    B35 = 'ARBITRUM'
    C35 = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
    D35 = '0.9515'
    E35 = '  +3.75%  '
    B36 = 'ImmutableX'
    C36 = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
    D36 = '0.6052'
    E36 = '  +5.21%  '
    B37 = 'TrustWalletToken'
    C37 = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
    D37 = '0.9490'
    E37 = '  +17.75%  '
    B38 = 'MXToken'
    C38 = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    D38 = '2.520'
    E38 = '  -0.57%  '
    B39 = 'PaxDollar'
    C39 = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
    D39 = '1.0000'
    E39 = '  +0.15%  '
    B40 = 'VeChain'
    C40 = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    D40 = '0.01477'
    E40 = '  -3.63%  '
    B41 = 'Quant'
    C41 = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
    D41 = '100.27'
    E41 = '  +3.16%  '
    B42 = 'RenderToken'
    C42 = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    D42 = '1.838'
    E42 = '  +4.52%  '
    B43 = 'TheSandbox'
    C43 = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
    D43 = '0.3743'
    E43 = '  +2.01%  '
    B44 = 'FraxShare'
    C44 = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    D44 = '4.923'
    E44 = '  +4.84%  '
    B45 = 'Algorand'
    C45 = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
    D45 = '0.1113'
    E45 = '  +2.55%  '
    B46 = 'Aptos'
    C46 = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
    D46 = '6.180'
    E46 = '  +3.30%  '
    B47 = 'Cronos'
    C47 = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
    D47 = '0.05264'
    E47 = '  +1.28%  '
    B48 = 'Elrond'
    C48 = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
    D48 = '29.72'
    E48 = '  +1.59%  '
    B49 = 'EnergySwap'
    C49 = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    D49 = '7.395'
    E49 = '  +3.59%  '
    B50 = 'USDD'
    C50 = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
    D50 = '1.002'
    E50 = '  +0.70%  '
    E51 = '  +0.27%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}

